$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Disaster" column header + values
$ws.Range("I1:I7").ClearFormats()
$ws.Range("I1").Value = "Disaster"
$ws.Range("I2").Value = 0.25
$ws.Range("I3").Value = 0.05
$ws.Range("I4").Value = 0.1
$ws.Range("I5").Value = 0.075
$ws.Range("I6").Value = 0.15
$ws.Range("I7").Value = 0.13

# Apply font color (black) to the new column cells
$ws.Range("I1:I7").Font.Color = 0

# Update selection to match final state
$ws.Range("I8").Select()
